# Atualiza datasets das ligas
# Insere o novo time "Pepe Leal FC" (id 1326835) na linha 15, empurrando
# os times subsequentes (Pontaç0 F.C., SC 100 Sono, SC ÉoINTER!, Texas Club 2026)
# uma linha para baixo, e recria os hyperlinks da coluna C na ordem correta.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Desloca as linhas 15:18 para 16:19, abrindo espaço para o novo time.
$ws.Rows("15:15").Insert()

# Preenche os dados do novo time na linha 15.
$ws.Range("A15").Value = "Pepe Leal FC"
$ws.Range("B15").Value = 1326835
$ws.Range("C15").Value = "https://cartola.globo.com/#!/time/1326835"
$ws.Range("C15").Style = "Hyperlink"

# Remove todos os hyperlinks existentes (o Insert() nao reindexa os refs)
# e os recria na ordem correta, linha a linha, a partir dos IDs atuais.
$ws.Range("C2").Hyperlinks.Delete()

$lastRow = $ws.Range("A1").End(-4121).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $teamId = $ws.Range("B" + $r).Value2
    $location = "!/time/" + $teamId
    $ws.Hyperlinks.Add($ws.Range("C" + $r), "https://cartola.globo.com/", $location) | Out-Null
    $ws.Range("C" + $r).Style = "Hyperlink"
}
